# Models_and_trailers.xlsx - "Tables changed to test possible conflict with
# multiple loadbuilders": grow the Trailers table with two new rows and flag
# a few Models rows as sourced from a new plant ("Juarez 2").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Grow the "Trailers" table from K1:S4 to K1:S6 (two new rows) -------
$trailers = $ws.ListObjects.Item("Trailers")
$trailers.Resize($ws.Range("K1:S6"))

# Row 5 is an interior table row -> copy the formatting used by the other
# interior rows (row 2). Row 6 becomes the new last row of the table -> copy
# the formatting that used to belong to the old last row (row 4).
$ws.Range("K2:S2").Copy()
$ws.Range("K5:S5").PasteSpecial(-4122)
$ws.Range("K4:S4").Copy()
$ws.Range("K6:S6").PasteSpecial(-4122)

# Rows with the larger (12pt) table font need the taller row height.
$ws.Rows.Item(5).RowHeight = 16
$ws.Rows.Item(6).RowHeight = 16

# New Trailers row: 10 FLATBED, Juarez -> Juarez 2, 576x102x120, overhang 1
$ws.Range("K5").Value = 10
$ws.Range("L5").Value = "FLATBED"
$ws.Range("M5").Value = "Juarez"
$ws.Range("N5").Value = "Juarez 2"
$ws.Range("O5").Value = 576
$ws.Range("P5").Value = 102
$ws.Range("Q5").Value = 120
$ws.Range("R5").Value = 1
$ws.Range("S5").Value = 1

# New Trailers row: 2 FLATBED, Juarez -> Juarez 2, 636x102x120, overhang 0
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = "FLATBED"
$ws.Range("M6").Value = "Juarez"
$ws.Range("N6").Value = "Juarez 2"
$ws.Range("O6").Value = 636
$ws.Range("P6").Value = 102
$ws.Range("Q6").Value = 120
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 1

# --- 2. Update a few Models rows ------------------------------------------
$ws.Range("A5").Value = 23
$ws.Range("C5").Value = "Juarez 2"

$ws.Range("A10").Value = 24

$ws.Range("A15").Value = 24
$ws.Range("C15").Value = "Juarez 2"

# --- 3. Move the active selection ------------------------------------------
$ws.Application.Goto($ws.Range("R10"))
